# Refactor: removed BillWiseSalesReport test case (TC004) row from the
# Pharmacy Reports SNCH test pack sheet. This deletes the entire row 4
# (Pharmacy\Reports\Sales\TC01BillWiseSalesReport.py.py / TC004), which
# shifts every subsequent row up by one and lets Excel recompute the
# worksheet dimension and shared-string table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row (not just clear contents) so remaining rows shift up.
$ws.Rows.Item(4).Delete()

# Update the saved selection to match the post-edit cursor position.
[void]$ws.Range("C20").Select()
